# Auto-generated edit script: Add data for 2025-09-19
# Applies cell-level numeric updates across multiple worksheets

$wb = $excel.ActiveWorkbook

# --- Citywide Totals ---
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('H2').Value = 84
$ws.Range('K2').Value = 107
$ws.Range('B3').Value = 62
$ws.Range('E3').Value = 102
$ws.Range('F3').Value = 98
$ws.Range('G3').Value = 97
$ws.Range('I3').Value = 152
$ws.Range('J3').Value = 157
$ws.Range('L3').Value = 179
$ws.Range('H5').Value = 4
$ws.Range('C6').Value = 357
$ws.Range('D6').Value = 315
$ws.Range('E6').Value = 328
$ws.Range('G6').Value = 364
$ws.Range('H6').Value = 337
$ws.Range('J6').Value = 301
$ws.Range('K6').Value = 384
$ws.Range('L6').Value = 349
$ws.Range('B7').Value = 388
$ws.Range('C7').Value = 480
$ws.Range('D7').Value = 492
$ws.Range('E7').Value = 490
$ws.Range('F7').Value = 564
$ws.Range('G7').Value = 531
$ws.Range('H7').Value = 537
$ws.Range('I7').Value = 658
$ws.Range('J7').Value = 563
$ws.Range('K7').Value = 683
$ws.Range('L7').Value = 649

# --- Grand Crossing ---
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('L6').Value = 35
$ws.Range('L7').Value = 48

# --- Armour Square ---
$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('D5').Value = 3
$ws.Range('D6').Value = 5

# --- Humboldt Park ---
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J3').Value = 5
$ws.Range('J6').Value = 12

# --- South Chicago ---
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('E3').Value = 5
$ws.Range('E5').Value = 8

# --- Englewood ---
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('B3').Value = 3
$ws.Range('I3').Value = 12
$ws.Range('G6').Value = 25
$ws.Range('B7').Value = 29
$ws.Range('G7').Value = 34
$ws.Range('I7').Value = 38

# --- By Neighborhood ---
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('D5').Value = 5
$ws.Range('E8').Value = 36
$ws.Range('G8').Value = 25
$ws.Range('H8').Value = 46
$ws.Range('J8').Value = 33
$ws.Range('K8').Value = 35
$ws.Range('L8').Value = 24
$ws.Range('G18').Value = 4
$ws.Range('B28').Value = 29
$ws.Range('G28').Value = 34
$ws.Range('I28').Value = 38
$ws.Range('L36').Value = 48
$ws.Range('J41').Value = 12
$ws.Range('K51').Value = 3
$ws.Range('D53').Value = 60
$ws.Range('F53').Value = 57
$ws.Range('H53').Value = 67
$ws.Range('C70').Value = 6
$ws.Range('K70').Value = 18
$ws.Range('D72').Value = 5
$ws.Range('L77').Value = 22
$ws.Range('E80').Value = 8
$ws.Range('B98').Value = 388
$ws.Range('C98').Value = 480
$ws.Range('D98').Value = 492
$ws.Range('E98').Value = 490
$ws.Range('F98').Value = 564
$ws.Range('G98').Value = 531
$ws.Range('H98').Value = 537
$ws.Range('I98').Value = 658
$ws.Range('J98').Value = 563
$ws.Range('K98').Value = 683
$ws.Range('L98').Value = 649

# --- Loop ---
$ws = $wb.Worksheets.Item('Loop')
$ws.Range('H2').Value = 8
$ws.Range('F3').Value = 10
$ws.Range('H5').Value = 1
$ws.Range('D6').Value = 34
$ws.Range('D7').Value = 60
$ws.Range('F7').Value = 57
$ws.Range('H7').Value = 67

# --- Printers Row ---
$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range('D4').Value = 5
$ws.Range('D5').Value = 5

# --- Roseland ---
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('L6').Value = 16
$ws.Range('L7').Value = 22

# --- Calumet Heights ---
$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('F3').Value = 1
$ws.Range('F5').Value = 4

# --- Little Village ---
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('I2').Value = 1
$ws.Range('I6').Value = 3

# --- Old Town ---
$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('C4').Value = 6
$ws.Range('K4').Value = 14
$ws.Range('C5').Value = 6
$ws.Range('K5').Value = 18

# --- Austin ---
$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L3').Value = 10
$ws.Range('E5').Value = 26
$ws.Range('G5').Value = 19
$ws.Range('H5').Value = 37
$ws.Range('J5').Value = 17
$ws.Range('K5').Value = 19
$ws.Range('E6').Value = 36
$ws.Range('G6').Value = 25
$ws.Range('H6').Value = 46
$ws.Range('J6').Value = 33
$ws.Range('K6').Value = 35
$ws.Range('L6').Value = 24
